# Adds the "Tussen Evaluatie" section (3 new paragraphs) right after the
# "Plaatjes van de kaart" paragraph, and relocates the hidden "_GoBack"
# bookmark into its own trailing empty paragraph, as in the target diff.

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the list: "Plaatjes van de kaart".
# It also carries the (hidden) _GoBack bookmark right after its text.
$targetText = "Plaatjes van de kaart"
$anchorParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq $targetText) {
        $anchorParagraph = $p
    }
}

if ($null -eq $anchorParagraph) {
    throw "Could not find paragraph containing '$targetText'"
}

$anchorRange = $anchorParagraph.Range
$insertPos = $anchorRange.Start + $targetText.Length
$insertPoint = $d.Range($insertPos, $insertPos)

# The _GoBack bookmark currently sits right after "kaart" (still inside this
# paragraph). Remove it here - it will be re-created further down, in its
# own new trailing paragraph, after the inserted text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert the three new paragraphs plus an extra (initially empty) paragraph
# that will hold the relocated _GoBack bookmark, all in one shot so the
# bookmark lands cleanly in its own paragraph rather than spanning a
# paragraph-mark boundary.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParagraphsXml =
  '<w:p ' + $wNs + '>' +
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Tussen Evaluatie</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wNs + '>' +
    '<w:r><w:t>Onze communicatie was niet optimaal, geen projectleider</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $wNs + '>' +
    '<w:r><w:t xml:space="preserve">Ongeorganiseerd werken met drive en </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>github</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>' +
  '<w:p ' + $wNs + '>' +
    '<w:bookmarkStart w:id="100" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="100"/>' +
  '</w:p>'

[void]$insertPoint.InsertXML($newParagraphsXml)
